# Add data for 2021-10-30 (updates the "through 10-21" snapshot to "through 10-22")
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet (tab) to reflect the new "as of" date
$ws.Name = "Through 2021-10-22"

# Update the row-12 label to match the new "as of" date
$ws.Range("A12").Value = "October (through 10-22)"

# --- Row 12 (October, partial-month row) updates ---
# 2016
$ws.Range("F12").Value = 34
$ws.Range("G12").Value = 0.08110000000000001
# 2017
$ws.Range("I12").Value = 40
$ws.Range("J12").Value = 0.1489
# 2018
$ws.Range("L12").Value = 48
$ws.Range("M12").Value = 0.0588
# 2019
$ws.Range("O12").Value = 32
$ws.Range("P12").Value = 0.1111
# 2020
$ws.Range("R12").Value = 103
# 2021
$ws.Range("U12").Value = 135

# --- Row 13 (Total) updates ---
# 2016
$ws.Range("F13").Value = 417
$ws.Range("G13").Value = 0.1052
# 2017
$ws.Range("I13").Value = 617
$ws.Range("J13").Value = 0.08459999999999999
# 2018
$ws.Range("L13").Value = 535
$ws.Range("M13").Value = 0.1068
# 2019
$ws.Range("O13").Value = 411
$ws.Range("P13").Value = 0.1026
# 2020
$ws.Range("R13").Value = 951
$ws.Range("S13").Value = 0.0528
# 2021
$ws.Range("U13").Value = 1300
$ws.Range("V13").Value = 0.0593
